$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.177215099334717
$ws.Range("B1").Value = 2.419754266738892
$ws.Range("D1").Value = 2.33219575881958
$ws.Range("E1").Value = 1.201651096343994
